$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (column D) and Volume(1h) (column E) cells with the
# latest crypto snapshot values. Most source values are plain text (prices
# use '.' as a thousands separator so some contain more than one '.', and
# percentages are padded with spaces and a trailing '%'), so for the few
# price values that otherwise look like plain decimal numbers we assign them
# through Formula with a leading apostrophe (the standard Excel 'text-quote'
# prefix) so they are kept as text instead of being parsed as numbers.

$ws.Range("D2").Value = '93.227.15'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '3.097.24'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Formula = "'236.69"
$ws.Range("E5").Value = '  -3.94%  '
$ws.Range("D6").Formula = "'612.45"
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("D7").Formula = "'1.13"
$ws.Range("E7").Value = '  +1.95%  '
$ws.Range("E8").Value = '  +0.98%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Formula = "'0.826"
$ws.Range("E10").Value = '  +12.24%  '
$ws.Range("D11").Value = '3.095.41'
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("E12").Value = '  -3.24%  '
$ws.Range("E13").Value = '  -3.20%  '
$ws.Range("D14").Formula = "'34.78"
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = '93.007.83'
$ws.Range("E16").Value = '  -3.46%  '
$ws.Range("D17").Value = '3.674.35'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '3.105.81'
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").Formula = "'3.67"
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("D20").Formula = "'14.65"
$ws.Range("E20").Value = '  -1.81%  '
$ws.Range("D21").Formula = "'5.97"
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").Formula = "'439.67"
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D23").Formula = "'0.0000199"
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("D24").Formula = "'9.00"
$ws.Range("E24").Value = '  -5.76%  '
$ws.Range("E25").Value = '  +4.61%  '
$ws.Range("D26").Formula = "'5.67"
$ws.Range("E26").Value = '  -3.51%  '
$ws.Range("D27").Formula = "'12.75"
$ws.Range("E27").Value = '  +8.22%  '
$ws.Range("D28").Formula = "'85.74"
$ws.Range("E28").Value = '  -2.57%  '
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").Formula = "'0.249"
$ws.Range("E30").Value = '  +4.94%  '
$ws.Range("E31").Value = '  +8.08%  '
$ws.Range("E32").Value = '  -15.37%  '
$ws.Range("E33").Value = '  -1.92%  '
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("D35").Formula = "'7.95"
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("E36").Value = '  -10.15%  '
$ws.Range("D37").Formula = "'25.85"
$ws.Range("E37").Value = '  -1.81%  '
$ws.Range("D38").Formula = "'3.96"
$ws.Range("E38").Value = '  -6.32%  '
$ws.Range("D39").Formula = "'1.89"
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("D40").Formula = "'0.448"
$ws.Range("E40").Value = '  +1.56%  '
$ws.Range("E41").Value = '  +8.06%  '
$ws.Range("D42").Formula = "'473.63"
$ws.Range("E42").Value = '  -3.74%  '
$ws.Range("D43").Formula = "'1.28"
$ws.Range("E43").Value = '  -2.09%  '
$ws.Range("D44").Formula = "'3.26"
$ws.Range("E44").Value = '  -4.29%  '
$ws.Range("D46").Formula = "'158.95"
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").Formula = "'0.698"
$ws.Range("E47").Value = '  -1.51%  '
$ws.Range("D48").Formula = "'1.85"
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").Formula = "'1.31"
$ws.Range("E49").Value = '  -3.56%  '
$ws.Range("D50").Formula = "'43.81"
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("D51").Formula = "'4.36"
$ws.Range("E51").Value = '  -1.57%  '
